# Update the dSF (column F) values for the data rows on Sheet1.
# These values reflect a re-pull/re-push of data with an updated
# mean calculation, changing several previously-zero placeholder
# values to their computed deviations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    4  = 2
    5  = 2
    6  = 1
    7  = -3
    8  = -2
    9  = 10
    10 = -2
    11 = -3
    12 = -4
    13 = 3
    14 = 2
    15 = -2
    16 = -2
    17 = 4
    18 = 4
    19 = -1
    20 = 3
    22 = 8
    23 = 1
    24 = 3
    25 = 2
    26 = 3
    27 = -3
    28 = 2
    29 = 2
    30 = -1
    31 = -2
    32 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
